# Auto update Excel log: append new mmWave sensor log rows (75-79)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-02-01", "11:38:33", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:38:42", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:38:53", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:39:16", "11:00", "Living Room", "NO_MOTION_DETECTED", "Inactive"),
    @("2026-02-01", "11:39:26", "11:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 75
$endRow = $startRow + $rows.Count - 1

# Column A holds date-shaped strings ("2026-02-01"). Writing them straight
# through .Value makes Excel auto-parse them into real date serials, so
# force the column-A cells to Text first, write the values, then restore
# the default "Normal" style so no extra number-format styling lingers on
# the cells.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

$dateRange.Style = "Normal"
